$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.861.51"
$ws.Range("E2").Value = "  +2.47%  "
$ws.Range("D3").Value = "3.088.64"
$ws.Range("E3").Value = "  +4.97%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.96"
$ws.Range("E5").Value = "  +2.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.58"
$ws.Range("E6").Value = "  +6.01%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.084.97"
$ws.Range("E8").Value = "  +4.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.524"
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.61"
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("E11").Value = "  +3.42%  "
$ws.Range("E12").Value = "  +5.76%  "
$ws.Range("E13").Value = "  +1.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.43"
$ws.Range("E14").Value = "  +5.97%  "
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "3.600.66"
$ws.Range("E16").Value = "  +4.92%  "
$ws.Range("D17").Value = "66.863.47"
$ws.Range("E17").Value = "  +2.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.19"
$ws.Range("E18").Value = "  +3.84%  "
$ws.Range("D19").Value = "3.090.01"
$ws.Range("E19").Value = "  +5.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.16"
$ws.Range("E20").Value = "  +9.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "466.28"
$ws.Range("E21").Value = "  +4.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.714"
$ws.Range("E22").Value = "  +4.05%  "
$ws.Range("E23").Value = "  +3.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.42"
$ws.Range("E24").Value = "  +1.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("E25").Value = "  +6.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.00"
$ws.Range("E26").Value = "  +7.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.14"
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.99"
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.39"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  +3.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000103"
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.25"
$ws.Range("E33").Value = "  +4.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.114"
$ws.Range("E34").Value = "  +3.48%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +3.26%  "
$ws.Range("E37").Value = "  +3.01%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.11"
$ws.Range("E38").Value = "  +6.85%  "
$ws.Range("B39").Value = "Arweave"
$ws.Range("C39").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "47.02"
$ws.Range("E39").Value = "  +5.93%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.26"
$ws.Range("E40").Value = "  +1.35%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.318"
$ws.Range("E41").Value = "  +6.72%  "
$ws.Range("E42").Value = "  +1.59%  "
$ws.Range("E43").Value = "  +2.72%  "
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0360"
$ws.Range("E45").Value = "  +2.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "384.69"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("D47").Value = "2.772.93"
$ws.Range("E47").Value = "  +2.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.12"
$ws.Range("E48").Value = "  +1.62%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("E50").Value = "  +6.50%  "
$ws.Range("E51").Value = "  +1.70%  "
